# Auto-generated Excel COM-interop script to apply price/profit updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21 (item id 2149) in ALC
$ws.Range("H21").Value = 3749
$ws.Range("I21").Value = 3498
$ws.Range("J21").Value = 4000
$ws.Range("K21").Value = 3498
$ws.Range("L21").Value = 4000
$ws.Range("M21").Value = -3030
$ws.Range("N21").Value = -4936

# Row 23 (item id 2149) in ALC
$ws.Range("H23").Value = 3749
$ws.Range("I23").Value = 3498
$ws.Range("J23").Value = 4000
$ws.Range("K23").Value = 3498
$ws.Range("L23").Value = 4000
$ws.Range("M23").Value = -3264
$ws.Range("N23").Value = -4468

# Row 33 (item id 5512) in ALC
$ws.Range("H33").Value = 370.26923
$ws.Range("I33").Value = 356.38095
$ws.Range("J33").Value = 428.6
$ws.Range("K33").Value = 356.38095
$ws.Range("L33").Value = 428.6
$ws.Range("M33").Value = -127.38095
$ws.Range("N33").Value = -886.6

# Row 40 (item id 5505) in ALC
$ws.Range("H40").Value = 16030.692
$ws.Range("J40").Value = 16533.25
$ws.Range("L40").Value = 16533.25
$ws.Range("N40").Value = -16883.25

# Row 92 (item id 19901) in ALC
$ws.Range("H92").Value = 431.7647
$ws.Range("I92").Value = 421.25
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 421.25
$ws.Range("L92").Value = 600
$ws.Range("M92").Value = 826.75
$ws.Range("N92").Value = -3096

# Row 113 (item id 27775) in ALC
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746

# Row 116 (item id 27778) in ALC
$ws.Range("H116").Value = 4176296.5
$ws.Range("I116").Value = 10303.333
$ws.Range("K116").Value = 10303.333
$ws.Range("M116").Value = -6861.333000000001

# Row 117 (item id 26118) in ALC
$ws.Range("H117").Value = 89290.11
$ws.Range("J117").Value = 89290.11
$ws.Range("L117").Value = 89290.11
$ws.Range("N117").Value = -98468.11

# Row 123 (item id 34090) in ALC
$ws.Range("H123").Value = 82457
$ws.Range("J123").Value = 82457
$ws.Range("L123").Value = 82457
$ws.Range("N123").Value = -92257

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (item id 27713) in ARM
$ws.Range("H2").Value = 1028.3914
$ws.Range("I2").Value = 913.63635
$ws.Range("K2").Value = 913.63635
$ws.Range("M2").Value = -800.63635

# Row 32 (item id 44147) in ARM
$ws.Range("H32").Value = 9580.02
$ws.Range("I32").Value = 4512.385
$ws.Range("J32").Value = 24782.924
$ws.Range("K32").Value = 4512.385
$ws.Range("L32").Value = 24782.924
$ws.Range("M32").Value = -4225.385
$ws.Range("N32").Value = -25356.924

# Row 45 (item id 27714) in ARM
$ws.Range("H45").Value = 2311.4666
$ws.Range("I45").Value = 1931.4166
$ws.Range("K45").Value = 1931.4166
$ws.Range("M45").Value = -1554.4166

# Row 63 (item id 12528) in ARM
$ws.Range("H63").Value = 2700.5715
$ws.Range("I63").Value = 2700.5715
$ws.Range("K63").Value = 2700.5715
$ws.Range("M63").Value = -2014.5715

# Row 66 (item id 12528) in ARM
$ws.Range("H66").Value = 2700.5715
$ws.Range("I66").Value = 2700.5715
$ws.Range("K66").Value = 13502.8575
$ws.Range("M66").Value = -10070.8575

# Row 97 (item id 19941) in ARM
$ws.Range("H97").Value = 621.55554
$ws.Range("I97").Value = 726.1429000000001
$ws.Range("J97").Value = 255.5
$ws.Range("K97").Value = 726.1429000000001
$ws.Range("L97").Value = 255.5
$ws.Range("M97").Value = -230.1429000000001
$ws.Range("N97").Value = -1247.5

# Row 116 (item id 27713) in ARM
$ws.Range("H116").Value = 1028.3914
$ws.Range("I116").Value = 913.63635
$ws.Range("K116").Value = 913.63635
$ws.Range("M116").Value = 1380.36365

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (item id 27713) in BSM
$ws.Range("H3").Value = 1028.3914
$ws.Range("I3").Value = 913.63635
$ws.Range("K3").Value = 913.63635
$ws.Range("M3").Value = -799.63635

# Row 107 (item id 27706) in BSM
$ws.Range("H107").Value = 2815.2273
$ws.Range("I107").Value = 2198.7334
$ws.Range("K107").Value = 2198.7334
$ws.Range("M107").Value = -278.7334000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (item id 27691) in CRP
$ws.Range("H16").Value = 1626.4615
$ws.Range("I16").Value = 1362.619
$ws.Range("J16").Value = 2734.6
$ws.Range("K16").Value = 1362.619
$ws.Range("L16").Value = 2734.6
$ws.Range("M16").Value = -1075.619
$ws.Range("N16").Value = -3308.6

# Row 113 (item id 27691) in CRP
$ws.Range("H113").Value = 1626.4615
$ws.Range("I113").Value = 1362.619
$ws.Range("J113").Value = 2734.6
$ws.Range("K113").Value = 1362.619
$ws.Range("L113").Value = 2734.6
$ws.Range("M113").Value = 807.3810000000001
$ws.Range("N113").Value = -7074.6

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (item id 4847) in CUL
$ws.Range("H2").Value = 123.125
$ws.Range("I2").Value = 158
$ws.Range("K2").Value = 948
$ws.Range("M2").Value = -835

# Row 3 (item id 44094) in CUL
$ws.Range("H3").Value = 3244.375
$ws.Range("I3").Value = 988.75
$ws.Range("J3").Value = 5500
$ws.Range("K3").Value = 2966.25
$ws.Range("L3").Value = 16500
$ws.Range("M3").Value = -2854.25
$ws.Range("N3").Value = -16724

# Row 103 (item id 19839) in CUL
$ws.Range("H103").Value = 667.55554
$ws.Range("I103").Value = 772.8570999999999
$ws.Range("K103").Value = 2318.5713
$ws.Range("M103").Value = -1439.5713

# Row 117 (item id 27870) in CUL
$ws.Range("H117").Value = 1408.2
$ws.Range("I117").Value = 241.66667
$ws.Range("K117").Value = 725.00001
$ws.Range("M117").Value = 2716.99999

# Row 131 (item id 36060) in CUL
$ws.Range("H131").Value = 34817.133
$ws.Range("J131").Value = 1739.762
$ws.Range("L131").Value = 5219.286
$ws.Range("N131").Value = -15299.286

# Row 141 (item id 44076) in CUL
$ws.Range("H141").Value = 2612.6924
$ws.Range("I141").Value = 2612.6924
$ws.Range("K141").Value = 7838.0772
$ws.Range("M141").Value = -2658.0772

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (item id 36182) in GSM
$ws.Range("H122").Value = 60491.105
$ws.Range("I122").Value = 70573.31
$ws.Range("K122").Value = 211719.93
$ws.Range("M122").Value = -209269.93

# Row 126 (item id 36184) in GSM
$ws.Range("H126").Value = 2699.2
$ws.Range("I126").Value = 2266.818
$ws.Range("K126").Value = 6800.454000000001
$ws.Range("M126").Value = -4330.454000000001

# Row 140 (item id 42458) in GSM
$ws.Range("H140").Value = 75308.82000000001
$ws.Range("J140").Value = 79049.625
$ws.Range("L140").Value = 79049.625
$ws.Range("N140").Value = -89409.625

# Row 141 (item id 42504) in GSM
$ws.Range("H141").Value = 157499.5
$ws.Range("J141").Value = 157499.5
$ws.Range("L141").Value = 157499.5
$ws.Range("N141").Value = -167859.5

$ws = $wb.Worksheets.Item("LTW")
# Row 51 (item id 3423) in LTW
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 61 (item id 27740) in LTW
$ws.Range("H61").Value = 1000.36365
$ws.Range("I61").Value = 900.4
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 900.4
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -698.4
$ws.Range("N61").Value = -2404

# Row 94 (item id 18067) in LTW
$ws.Range("H94").Value = 29665.5
$ws.Range("J94").Value = 29665.5
$ws.Range("L94").Value = 29665.5
$ws.Range("N94").Value = -31017.5

# Row 113 (item id 27740) in LTW
$ws.Range("H113").Value = 1000.36365
$ws.Range("I113").Value = 900.4
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 900.4
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1269.6
$ws.Range("N113").Value = -6340

# Row 122 (item id 36247) in LTW
$ws.Range("H122").Value = 10021418
$ws.Range("I122").Value = 35254.816
$ws.Range("K122").Value = 105764.448
$ws.Range("M122").Value = -103314.448

# Row 127 (item id 34401) in LTW
$ws.Range("H127").Value = 59926.668
$ws.Range("J127").Value = 59926.668
$ws.Range("L127").Value = 59926.668
$ws.Range("N127").Value = -69846.66800000001

# Row 132 (item id 44058) in LTW
$ws.Range("H132").Value = 14548.363
$ws.Range("I132").Value = 16355.421
$ws.Range("J132").Value = 3103.6667
$ws.Range("K132").Value = 49066.263
$ws.Range("L132").Value = 9311.000100000001
$ws.Range("M132").Value = -46536.263
$ws.Range("N132").Value = -14371.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (item id 19981) in WVR
$ws.Range("H100").Value = 2646965.2
$ws.Range("I100").Value = 4763597
$ws.Range("J100").Value = 1175.8334
$ws.Range("K100").Value = 9527194
$ws.Range("L100").Value = 2351.6668
$ws.Range("M100").Value = -9526653
$ws.Range("N100").Value = -3433.6668

# Row 107 (item id 27746) in WVR
$ws.Range("H107").Value = 1828.8
$ws.Range("I107").Value = 1375.0714
$ws.Range("J107").Value = 2073.1155
$ws.Range("K107").Value = 4125.2142
$ws.Range("L107").Value = 6219.3465
$ws.Range("M107").Value = -2205.2142
$ws.Range("N107").Value = -10059.3465

# Row 132 (item id 44029) in WVR
$ws.Range("H132").Value = 791601.9399999999
$ws.Range("J132").Value = 2900079.2
$ws.Range("L132").Value = 8700237.600000001
$ws.Range("N132").Value = -8705297.600000001

# Row 135 (item id 42043) in WVR
$ws.Range("H135").Value = 83995
$ws.Range("J135").Value = 83995
$ws.Range("L135").Value = 83995
$ws.Range("N135").Value = -94135

# Row 136 (item id 44031) in WVR
$ws.Range("H136").Value = 2656.1292
$ws.Range("I136").Value = 2502.111
$ws.Range("K136").Value = 7506.333
$ws.Range("M136").Value = -4956.333
